$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 229, pushing existing rows 229-288
# down to 230-289 (dimension grows from A1:R288 to A1:R289).
$ws.Rows.Item(229).Insert()

# Populate the newly inserted row 229 with the new record.
$ws.Range("A229").Value = 1
$ws.Range("B229").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C229").Value = "Arica y Parinacota"
$ws.Range("D229").Value = 44722
$ws.Range("E229").Value = 15
$ws.Range("F229").Value = 100114013
$ws.Range("G229").Value = "Zanahoria"
$ws.Range("H229").Value = "Sin especificar"
$ws.Range("I229").Value = "Primera"
$ws.Range("J229").Value = 80
$ws.Range("K229").Value = 14000
$ws.Range("L229").Value = 15000
$ws.Range("M229").Value = 14500
$ws.Range("N229").Value = "`$/saco 25 kilos"
$ws.Range("O229").Value = "Valle de Camiña"
$ws.Range("P229").Value = 580
$ws.Range("Q229").Value = 25
$ws.Range("R229").Value = "Hortaliza"
